$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.039.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.888.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4737"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3951"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.07"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08025"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.020"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.885.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.031"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.218"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.017"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001053"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.021.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.522"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.350"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.109.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.112"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.519"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9755"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09585"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.644"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.353"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.368"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06084"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02251"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.201"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.195"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.013"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5964"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1892"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.269"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5658"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.930"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.371"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06827"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.90%  "
